$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value()
    if ($val -eq $oldValue) {
        $cell.Value = $newValue
    }
}
